$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 2")
$ws.Activate()

# Fill in new log entries for rows 26-29
$ws.Range("A26").Value = "Working on getting InfluxDB and Telegraf running locally"
$ws.Range("B26").Value = "Parsing Files"
$ws.Range("C26").Value = 1.25

$ws.Range("A27").Value = "Logs and GitHub cleanup"
$ws.Range("B27").Value = "Project Management"
$ws.Range("C27").Value = 0.25

$ws.Range("A28").Value = "Research how streaming anomaly detection would wrok with InfluxDB"
$ws.Range("B28").Value = "InfluxDB Streaming"
$ws.Range("C28").Value = 0.25

$ws.Range("A29").Value = "Stand-up and presentation run-through"
$ws.Range("B29").Value = "Internal Meetings"
$ws.Range("C29").Value = 0.75

# Update selection to match final state
$ws.Range("A30").Select()
